{"js": "// Remove the \"Unsubscribe\" hyperlink (and its text run) from the document,\n// leaving the paragraph itself (and its indentation formatting) intact.\nconst body = context.document.body;\n\nconst results = body.search(\"Unsubscribe\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (const item of results.items) {\n  item.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Unsubscribe\" hyperlink text and remove it (including its\n# run), leaving the paragraph (and its indentation formatting) in place \u2014\n# mirrors deleting the <w:hyperlink> wrapper + run from the OOXML while\n# keeping the now-empty <w:p>.\n$rng = $d.Content\n$rng.Find.Text = \"Unsubscribe\"\n$rng.Find.Forward = $true\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $true\n$found = $rng.Find.Execute()\n\nif ($found) {\n    $rng.Delete()\n}\n"}
